$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 105, pushing the existing rows 105-132 down to 107-134.
$ws.Rows.Item(105).Resize(2).Insert()

# Populate the first new row (105) - Red Globe, week of 2022-04-12.
$ws.Cells.Item(105, 1).Value = 11
$ws.Cells.Item(105, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(105, 3).Value = "Bíobío"
$ws.Cells.Item(105, 4).Value = 44663
$ws.Cells.Item(105, 5).Value = 8
$ws.Cells.Item(105, 6).Value = "Fruta"
$ws.Cells.Item(105, 7).Value = 100109
$ws.Cells.Item(105, 8).Value = "Uva"
$ws.Cells.Item(105, 9).Value = 100109001
$ws.Cells.Item(105, 10).Value = "Uva"
$ws.Cells.Item(105, 11).Value = "Red Globe"
$ws.Cells.Item(105, 12).Value = "Primera"
$ws.Cells.Item(105, 13).Value = 160
$ws.Cells.Item(105, 14).Value = 9000
$ws.Cells.Item(105, 15).Value = 10000
$ws.Cells.Item(105, 16).Value = 9500
$ws.Cells.Item(105, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(105, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(105, 19).Value = 528
$ws.Cells.Item(105, 20).Value = 18

# Populate the second new row (106) - Thompson seedless, week of 2022-04-12.
$ws.Cells.Item(106, 1).Value = 11
$ws.Cells.Item(106, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(106, 3).Value = "Bíobío"
$ws.Cells.Item(106, 4).Value = 44663
$ws.Cells.Item(106, 5).Value = 8
$ws.Cells.Item(106, 6).Value = "Fruta"
$ws.Cells.Item(106, 7).Value = 100109
$ws.Cells.Item(106, 8).Value = "Uva"
$ws.Cells.Item(106, 9).Value = 100109001
$ws.Cells.Item(106, 10).Value = "Uva"
$ws.Cells.Item(106, 11).Value = "Thompson seedless"
$ws.Cells.Item(106, 12).Value = "Primera"
$ws.Cells.Item(106, 13).Value = 150
$ws.Cells.Item(106, 14).Value = 12000
$ws.Cells.Item(106, 15).Value = 13000
$ws.Cells.Item(106, 16).Value = 12533
$ws.Cells.Item(106, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(106, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(106, 19).Value = 696
$ws.Cells.Item(106, 20).Value = 18
